$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row, computed from a permutation of the existing rows
# (the D/H/I/J/K/L/M/N/O/P columns were reshuffled across rows 2-19;
# columns A,B,C,E,F,G,Q,R are identical in every row so are left untouched)

# Row 2 (was row 5)
$ws.Range("D2").Value = 44489
$ws.Range("H2").Value = 'Sin especificar'
$ws.Range("I2").Value = 'Primera'
$ws.Range("J2").Value = 600
$ws.Range("K2").Value = 1400
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = 1450
$ws.Range("N2").Value = '$/kilo'
$ws.Range("O2").Value = 'Provincia de Linares'
$ws.Range("P2").Value = 1450

# Row 3 (was row 7)
$ws.Range("D3").Value = 44511
$ws.Range("H3").Value = 'Sin especificar'
$ws.Range("I3").Value = 'Primera'
$ws.Range("J3").Value = 600
$ws.Range("K3").Value = 1300
$ws.Range("L3").Value = 1400
$ws.Range("M3").Value = 1350
$ws.Range("N3").Value = '$/kilo'
$ws.Range("O3").Value = 'Provincia de Linares'
$ws.Range("P3").Value = 1350

# Row 4 (was row 8)
$ws.Range("D4").Value = 44526
$ws.Range("H4").Value = 'Sin especificar'
$ws.Range("I4").Value = 'Primera'
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 1500
$ws.Range("L4").Value = 1600
$ws.Range("M4").Value = 1550
$ws.Range("N4").Value = '$/kilo'
$ws.Range("O4").Value = 'Provincia de Linares'
$ws.Range("P4").Value = 1550

# Row 5 (was row 19)
$ws.Range("D5").Value = 44468
$ws.Range("H5").Value = 'Verde'
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 1800
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = 1920
$ws.Range("N5").Value = '$/kilo'
$ws.Range("O5").Value = 'Provincia de Linares'
$ws.Range("P5").Value = 1920

# Row 6 (was row 13)
$ws.Range("D6").Value = 44519
$ws.Range("H6").Value = 'Sin especificar'
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 250
$ws.Range("K6").Value = 1200
$ws.Range("L6").Value = 1300
$ws.Range("M6").Value = 1240
$ws.Range("N6").Value = '$/kilo'
$ws.Range("O6").Value = 'Provincia de Linares'
$ws.Range("P6").Value = 1240

# Row 7 (was row 2)
$ws.Range("D7").Value = 44839
$ws.Range("H7").Value = 'Sin especificar'
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 500
$ws.Range("K7").Value = 1700
$ws.Range("L7").Value = 1800
$ws.Range("M7").Value = 1760
$ws.Range("N7").Value = '$/kilo'
$ws.Range("O7").Value = 'Provincia de Linares'
$ws.Range("P7").Value = 1760

# Row 8 (was row 9)
$ws.Range("D8").Value = 44876
$ws.Range("H8").Value = 'Sin especificar'
$ws.Range("I8").Value = 'Primera'
$ws.Range("J8").Value = 350
$ws.Range("K8").Value = 1500
$ws.Range("L8").Value = 1600
$ws.Range("M8").Value = 1557
$ws.Range("N8").Value = '$/kilo'
$ws.Range("O8").Value = 'Provincia de Linares'
$ws.Range("P8").Value = 1557

# Row 9 (was row 3)
$ws.Range("D9").Value = 44868
$ws.Range("H9").Value = 'Sin especificar'
$ws.Range("I9").Value = 'Primera'
$ws.Range("J9").Value = 1000
$ws.Range("K9").Value = 1200
$ws.Range("L9").Value = 1300
$ws.Range("M9").Value = 1250
$ws.Range("N9").Value = '$/kilo'
$ws.Range("O9").Value = 'Región del Maule'
$ws.Range("P9").Value = 1250

# Row 10 (was row 4)
$ws.Range("D10").Value = 44868
$ws.Range("H10").Value = 'Sin especificar'
$ws.Range("I10").Value = 'Segunda'
$ws.Range("J10").Value = 200
$ws.Range("K10").Value = 1000
$ws.Range("L10").Value = 1000
$ws.Range("M10").Value = 1000
$ws.Range("N10").Value = '$/kilo'
$ws.Range("O10").Value = 'Región del Maule'
$ws.Range("P10").Value = 1000

# Row 11 (was row 18)
$ws.Range("D11").Value = 44510
$ws.Range("H11").Value = 'Sin especificar'
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 600
$ws.Range("K11").Value = 1300
$ws.Range("L11").Value = 1400
$ws.Range("M11").Value = 1350
$ws.Range("N11").Value = '$/kilo'
$ws.Range("O11").Value = 'Provincia de Linares'
$ws.Range("P11").Value = 1350

# Row 12 (was row 6)
$ws.Range("D12").Value = 44860
$ws.Range("H12").Value = 'Sin especificar'
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 1100
$ws.Range("K12").Value = 1500
$ws.Range("L12").Value = 1700
$ws.Range("M12").Value = 1609
$ws.Range("N12").Value = '$/kilo'
$ws.Range("O12").Value = 'Provincia de Linares'
$ws.Range("P12").Value = 1609

# Row 13 (was row 14)
$ws.Range("D13").Value = 44545
$ws.Range("H13").Value = 'Sin especificar'
$ws.Range("I13").Value = 'Primera'
$ws.Range("J13").Value = 550
$ws.Range("K13").Value = 1700
$ws.Range("L13").Value = 1800
$ws.Range("M13").Value = 1755
$ws.Range("N13").Value = '$/kilo'
$ws.Range("O13").Value = 'Provincia de Linares'
$ws.Range("P13").Value = 1755

# Row 14 (was row 10)
$ws.Range("D14").Value = 44524
$ws.Range("H14").Value = 'Sin especificar'
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 200
$ws.Range("K14").Value = 1500
$ws.Range("L14").Value = 1600
$ws.Range("M14").Value = 1550
$ws.Range("N14").Value = '$/kilo'
$ws.Range("O14").Value = 'Provincia de Talca'
$ws.Range("P14").Value = 1550

# Row 15 (was row 12)
$ws.Range("D15").Value = 44477
$ws.Range("H15").Value = 'Sin especificar'
$ws.Range("I15").Value = 'Primera'
$ws.Range("J15").Value = 500
$ws.Range("K15").Value = 1400
$ws.Range("L15").Value = 1500
$ws.Range("M15").Value = 1460
$ws.Range("N15").Value = '$/kilo'
$ws.Range("O15").Value = 'Provincia de Linares'
$ws.Range("P15").Value = 1460

# Row 16 (was row 11)
$ws.Range("D16").Value = 44875
$ws.Range("H16").Value = 'Sin especificar'
$ws.Range("I16").Value = 'Primera'
$ws.Range("J16").Value = 300
$ws.Range("K16").Value = 1500
$ws.Range("L16").Value = 1600
$ws.Range("M16").Value = 1550
$ws.Range("N16").Value = '$/kilo'
$ws.Range("O16").Value = 'Provincia de Linares'
$ws.Range("P16").Value = 1550

# Row 17 (was row 15)
$ws.Range("D17").Value = 44881
$ws.Range("H17").Value = 'Sin especificar'
$ws.Range("I17").Value = 'Primera'
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 2600
$ws.Range("L17").Value = 2700
$ws.Range("M17").Value = 2650
$ws.Range("N17").Value = '$/kilo'
$ws.Range("O17").Value = 'Provincia de Linares'
$ws.Range("P17").Value = 2650

# Row 18 (was row 16)
$ws.Range("D18").Value = 44881
$ws.Range("H18").Value = 'Sin especificar'
$ws.Range("I18").Value = 'Segunda'
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 2400
$ws.Range("L18").Value = 2400
$ws.Range("M18").Value = 2400
$ws.Range("N18").Value = '$/kilo'
$ws.Range("O18").Value = 'Provincia de Linares'
$ws.Range("P18").Value = 2400

# Row 19 (was row 17)
$ws.Range("D19").Value = 44496
$ws.Range("H19").Value = 'Sin especificar'
$ws.Range("I19").Value = 'Primera'
$ws.Range("J19").Value = 550
$ws.Range("K19").Value = 1500
$ws.Range("L19").Value = 2000
$ws.Range("M19").Value = 1773
$ws.Range("N19").Value = '$/paquete'
$ws.Range("O19").Value = 'Provincia de Linares'
$ws.Range("P19").Value = 1773

